# Update "想去人数" (people interested) counts in column F for a handful of
# rows on both the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    19 = 217
    24 = 53
    27 = 724
    28 = 2524
    29 = 18
    31 = 505
    32 = 821
    33 = 564
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
